$d = $word.ActiveDocument

$pairs = @(
    @("670÷8=83, 6", "839÷7=119, 6"),
    @("687÷2=343, 1", "992÷6=165, 2"),
    @("394÷4=98, 2", "109÷5=21, 4"),
    @("179÷7=25, 4", "487÷7=69, 4"),
    @("190÷7=27, 1", "814÷6=135, 4"),
    @("102÷6=17, 0", "214÷9=23, 7"),
    @("193÷4=48, 1", "575÷3=191, 2"),
    @("955÷9=106, 1", "554÷5=110, 4"),
    @("329÷2=164, 1", "405÷4=101, 1"),
    @("702÷8=87, 6", "253÷6=42, 1"),
    @("171÷3=57, 0", "716÷5=143, 1"),
    @("947÷7=135, 2", "370÷4=92, 2"),
    @("461÷4=115, 1", "915÷2=457, 1"),
    @("879÷9=97, 6", "986÷2=493, 0"),
    @("935÷9=103, 8", "844÷5=168, 4"),
    @("180÷5=36, 0", "868÷4=217, 0"),
    @("782÷4=195, 2", "208÷3=69, 1"),
    @("897÷4=224, 1", "763÷2=381, 1"),
    @("961÷2=480, 1", "600÷4=150, 0"),
    @("640÷5=128, 0", "765÷9=85, 0"),
    @("615÷9=68, 3", "197÷4=49, 1"),
    @("667÷6=111, 1", "695÷2=347, 1"),
    @("616÷6=102, 4", "713÷8=89, 1"),
    @("945÷3=315, 0", "886÷8=110, 6"),
    @("793÷9=88, 1", "680÷7=97, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
